$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Replace spaces with underscores in the category labels
$ws.Range("A4").Value = "United_States"
$ws.Range("A5").Value = "Other_countries"

# Update the selected cell/range on the sheet to match the saved view state
$ws.Range("B12").Select()
